$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("WhenNullIsPassed", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
while ($found) {
    $s = $rng.Start

    # remove the old "WhenNullIsPassed" text
    $rng.Text = ""

    # re-insert the suffix first, then "If" in front of it, so the final
    # order reads "If" + "NullIsPassed" -- each InsertAfter() call produces
    # its own run.
    $r1 = $d.Range($s, $s)
    $r1.InsertAfter("NullIsPassed")

    $r2 = $d.Range($s, $s)
    $r2.InsertAfter("If")

    # toggling a character property on the "If" run and back forces the
    # engine to keep it as its own run instead of re-coalescing it with
    # its identically-formatted neighbours.
    $ifRng = $d.Range($s, $s + 2)
    $ifRng.Bold = 1
    $ifRng.Bold = 0

    $rng = $d.Content
    $rng.Start = $s + 15
    $found = $rng.Find.Execute("WhenNullIsPassed", $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
}
